# Update countries & provincias Spain
# - Pakistan's case counts were refreshed, moving it above Portugal in the
#   (descending-sorted) country list.
# - Belice and Nueva Caledonia swapped places in the sorted list.
# - A couple of other countries (Nepal, Mongolia) got small count updates.
# - The "last updated" timestamp footer was bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 06:04"

# --- Pakistan moves above Portugal with refreshed numbers --------------
$ws.Range("A24").Value = "Pakistan"
$ws.Range("B24").Value = 27474
$ws.Range("C24").Value = 1039
$ws.Range("D24").Value = 7756
$ws.Range("E24").Value = 19100
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 19
$ws.Range("H24").Value = 618

$ws.Range("A25").Value = "Portugal"
$ws.Range("B25").Value = 27268
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 2422
$ws.Range("E25").Value = 23732
$ws.Range("F25").Value = 127
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 1114

# --- Nepal / Mongolia small updates -------------------------------------
$ws.Range("D160").Value = 30
$ws.Range("E160").Value = 79

$ws.Range("D177").Value = 14
$ws.Range("E177").Value = 28

# --- Belice moves above Nueva Caledonia (unchanged data, just reordered) -
$ws.Range("A192").Value = "Belice"
$ws.Range("B192").Value = 18
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 16
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 2

$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 18
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0
